# Updated cryptos list on Thu Aug 17 20:43:10 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row with
# new snapshot values, and fixes the BabyDogeCoin / RenderToken rows which
# had swapped places (row 47 is now BabyDogeCoin, row 48 is RenderToken).
#
# Price values that look like a plain decimal number (e.g. "226.01") are
# written with a leading apostrophe so Excel keeps them as text instead of
# auto-converting them to a number - matching how these cells were already
# stored (as literal strings, not numerics) before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.826.01"
$ws.Range("E2").Value = "  -4.49%  "
$ws.Range("D3").Value = "1.732.96"
$ws.Range("E3").Value = "  -4.54%  "
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'226.01"
$ws.Range("E5").Value = "  -3.39%  "
$ws.Range("D6").Value = "'0.5740"
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("D8").Value = "'0.2724"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'22.91"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").Value = "'0.06595"
$ws.Range("D11").Value = "'0.07534"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "1.738.22"
$ws.Range("E12").Value = "  -4.40%  "
$ws.Range("D13").Value = "'4.686"
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "'0.5987"
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("D15").Value = "1.972.50"
$ws.Range("E15").Value = "  -4.33%  "
$ws.Range("D16").Value = "'74.21"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "'0.000008633"
$ws.Range("E17").Value = "  -9.61%  "
$ws.Range("D18").Value = "27.832.72"
$ws.Range("E18").Value = "  -3.74%  "
$ws.Range("D19").Value = "'5.288"
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'204.35"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").Value = "'11.22"
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").Value = "'6.586"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "'149.79"
$ws.Range("D26").Value = "'8.011"
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("D27").Value = "'0.1225"
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("D28").Value = "'16.14"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").Value = "'1.377"
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("D30").Value = "'0.06144"
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("D31").Value = "'1.389"
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("D32").Value = "'3.726"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "'3.709"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'1.668"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("D35").Value = "'1.029"
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("D36").Value = "'0.6415"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").Value = "'2.423"
$ws.Range("E37").Value = "  -4.76%  "
$ws.Range("D38").Value = "'2.696"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("D39").Value = "'0.01664"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").Value = "1.118.38"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").Value = "'6.152"
$ws.Range("E41").Value = "  -4.67%  "
$ws.Range("D42").Value = "'0.8705"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D44").Value = "'99.72"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("E45").Value = "  -4.48%  "
$ws.Range("D46").Value = "'59.13"
$ws.Range("E46").Value = "  -2.68%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000108"
$ws.Range("E47").Value = "  -5.14%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.561"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D49").Value = "'8.237"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("E51").Value = "  -2.90%  "
